# Data Inquiry.xlsx - "add number and fix bug edit sales"
#
# Two new inquiry rows were added to the table (ids 18732 and 18735),
# pushing the existing rows down. Net effect on the sheet grid:
#   old row 4 -> new row 5
#   old row 5 -> new row 6
#   (new row 4 inserted, brand new data)
#   (new row 7 inserted, brand new data)
#   old row 6 -> new row 8
#   old row 7 -> new row 9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A scratch cell used to push date-shaped strings ("2021-09-07", "2021-09-07
# 13:12:27", ...) into the grid as plain text instead of letting Excel's
# input parser silently convert them into date serial numbers. We build the
# text via a formula (so it is never re-parsed as a literal typed value),
# copy the computed result, and paste-special just the value into the real
# target cell - this keeps the destination cell a plain shared string with
# no number formatting applied, same as a normal text value.
function Set-TextCell($cellRef, [string]$text) {
    $scratch = $ws.Range("ZZ1")
    $escaped = $text -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.ClearContents()
}

# ---- insert the two new rows, shifting rows 4-7 down as needed ----
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(7).Insert()

# ---- new row 4 ----
$ws.Range("A4").Value = 18732
$ws.Range("B4").Value = "IT "
Set-TextCell "C4" "2021-09-07 13:12:27"
$ws.Range("D4").Value = "EATON"
$ws.Range("E4").Value = "eaton serie xxxxxx"
$ws.Range("F4").Value = 1
Set-TextCell "G4" "2021-09-07"
$ws.Range("H4").Value = "yeyeye"
$ws.Range("I4").Value = "STOCK"
$ws.Range("J4").Value = 1
Set-TextCell "K4" "2021-09-07 13:20:49"
$ws.Range("L4").Value = "test"
$ws.Range("M4").Value = 0.003
$ws.Range("N4").Value = "sgd"
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 100000
$ws.Range("Q4").Value = 100000
$ws.Range("R4").Value = 100000
$ws.Range("S4").Value = "ok"
$ws.Range("T4").Value = "MAYENTI "

# ---- new row 7 ----
$ws.Range("A7").Value = 18735
$ws.Range("B7").Value = "IT"
Set-TextCell "C7" "2021-09-03 11:17:24"
$ws.Range("D7").Value = "NACOL"
$ws.Range("E7").Value = "tttt"
$ws.Range("F7").Value = 1
Set-TextCell "G7" "2021-09-03"
$ws.Range("H7").Value = "tested"
$ws.Range("I7").Value = "PRICE+LT"
$ws.Range("J7").Value = 1
Set-TextCell "K7" "2021-09-06 16:46:12"
$ws.Range("L7").Value = "tttt"
$ws.Range("M7").Value = 23223.131
$ws.Range("N7").Value = "EURO"
$ws.Range("O7").Value = 12324
$ws.Range("P7").Value = 2232323
$ws.Range("Q7").Value = 1212121
$ws.Range("R7").Value = 121212
$ws.Range("S7").Value = "23 WEEKS"
$ws.Range("T7").Value = "IT "
